$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$excel.Goto($ws.Range("AD1"), $false)
Write-Host "VisibleRange: $($win.VisibleRange.Address())"
